$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.838.61"
$ws.Range("E2").Value = "  -5.44%  "
$ws.Range("D3").Value = "3.364.66"
$ws.Range("E3").Value = "  -6.72%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.08"
$ws.Range("E5").Value = "  -6.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.76"
$ws.Range("E6").Value = "  -8.56%  "
$ws.Range("E7").Value = "  -5.08%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "3.355.29"
$ws.Range("E9").Value = "  -6.68%  "
$ws.Range("E10").Value = "  -12.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.597"
$ws.Range("E11").Value = "  -7.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.96"
$ws.Range("E12").Value = "  -11.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -11.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.76"
$ws.Range("E14").Value = "  -9.31%  "
$ws.Range("D15").Value = "3.901.58"
$ws.Range("E15").Value = "  -6.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "604.25"
$ws.Range("E16").Value = "  -11.25%  "
$ws.Range("D17").Value = "66.759.23"
$ws.Range("E17").Value = "  -5.68%  "
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("D19").Value = "3.366.14"
$ws.Range("E19").Value = "  -6.41%  "
$ws.Range("E20").Value = "  -7.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.70"
$ws.Range("E21").Value = "  -8.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.916"
$ws.Range("E22").Value = "  -8.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.26"
$ws.Range("E23").Value = "  -8.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.09"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.10"
$ws.Range("E25").Value = "  -13.91%  "
$ws.Range("E26").Value = "  -10.53%  "
$ws.Range("E27").Value = "  -9.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.65"
$ws.Range("E28").Value = "  -9.40%  "
$ws.Range("E29").Value = "  -12.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.87"
$ws.Range("E30").Value = "  -10.40%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.38"
$ws.Range("E31").Value = "  -11.31%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.86"
$ws.Range("E32").Value = "  -14.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.27"
$ws.Range("E33").Value = "  -8.50%  "
$ws.Range("E34").Value = "  -7.97%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.66"
$ws.Range("E35").Value = "  -7.66%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "537.99"
$ws.Range("E36").Value = "  +4.30%  "
$ws.Range("D37").Value = "3.766.87"
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.80"
$ws.Range("E39").Value = "  +39.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"
$ws.Range("E40").Value = "  -4.55%  "
$ws.Range("E41").Value = "  -13.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").Value = "  -9.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.354"
$ws.Range("E43").Value = "  -8.16%  "
$ws.Range("E44").Value = "  -7.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.94"
$ws.Range("E45").Value = "  -10.67%  "
$ws.Range("E46").Value = "  -10.16%  "
$ws.Range("E47").Value = "  -6.98%  "
$ws.Range("E48").Value = "  -12.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("E49").Value = "  -8.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.77"
$ws.Range("E51").Value = "  -10.03%  "
